$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 8500
$ws.Range("J7").Value = 8500
$ws.Range("L7").Value = 8500
$ws.Range("N7").Value = -8724
$ws.Range("H14").Value = 8500
$ws.Range("J14").Value = 8500
$ws.Range("L14").Value = 8500
$ws.Range("N14").Value = -8882
$ws.Range("H41").Value = 998.08
$ws.Range("I41").Value = 1350.6666
$ws.Range("J41").Value = 672.61536
$ws.Range("K41").Value = 1350.6666
$ws.Range("L41").Value = 672.61536
$ws.Range("M41").Value = -910.6666
$ws.Range("N41").Value = -1552.61536
$ws.Range("H68").Value = 35200
$ws.Range("J68").Value = 35200
$ws.Range("L68").Value = 35200
$ws.Range("N68").Value = -36698
$ws.Range("H71").Value = 35200
$ws.Range("J71").Value = 35200
$ws.Range("L71").Value = 105600
$ws.Range("N71").Value = -113088
$ws.Range("H76").Value = 4858
$ws.Range("I76").Value = 4751.75
$ws.Range("J76").Value = 4999.6665
$ws.Range("K76").Value = 4751.75
$ws.Range("L76").Value = 4999.6665
$ws.Range("M76").Value = -4436.75
$ws.Range("N76").Value = -5629.6665
$ws.Range("H79").Value = 4858
$ws.Range("I79").Value = 4751.75
$ws.Range("J79").Value = 4999.6665
$ws.Range("K79").Value = 4751.75
$ws.Range("L79").Value = 4999.6665
$ws.Range("M79").Value = -3659.75
$ws.Range("N79").Value = -7183.6665
$ws.Range("H112").Value = 1378.8276
$ws.Range("J112").Value = 1410.2142
$ws.Range("L112").Value = 4230.642599999999
$ws.Range("N112").Value = -6446.642599999999
$ws.Range("H125").Value = 4540.3335
$ws.Range("I125").Value = 3220.8
$ws.Range("J125").Value = 5482.857
$ws.Range("K125").Value = 28987.2
$ws.Range("L125").Value = 49345.713
$ws.Range("M125").Value = -26527.2
$ws.Range("N125").Value = -54265.713
$ws.Range("H132").Value = 5562326.5
$ws.Range("I132").Value = 6586203
$ws.Range("J132").Value = 4140.7144
$ws.Range("K132").Value = 19758609
$ws.Range("L132").Value = 12422.1432
$ws.Range("M132").Value = -19756079
$ws.Range("N132").Value = -17482.1432
$ws.Range("H138").Value = 3442.06
$ws.Range("I138").Value = 3714.2856
$ws.Range("J138").Value = 3397.7441
$ws.Range("K138").Value = 11142.8568
$ws.Range("L138").Value = 10193.2323
$ws.Range("M138").Value = -6002.856800000001
$ws.Range("N138").Value = -20473.2323

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26105.293
$ws.Range("I32").Value = 7274.0547
$ws.Range("J32").Value = 178847.56
$ws.Range("K32").Value = 7274.0547
$ws.Range("L32").Value = 178847.56
$ws.Range("M32").Value = -6987.0547
$ws.Range("N32").Value = -179421.56
$ws.Range("H74").Value = 942.3200000000001
$ws.Range("I74").Value = 849.82355
$ws.Range("J74").Value = 1138.875
$ws.Range("K74").Value = 849.82355
$ws.Range("L74").Value = 1138.875
$ws.Range("M74").Value = 24.17645000000005
$ws.Range("N74").Value = -2886.875
$ws.Range("H77").Value = 942.3200000000001
$ws.Range("I77").Value = 849.82355
$ws.Range("J77").Value = 1138.875
$ws.Range("K77").Value = 4249.117749999999
$ws.Range("L77").Value = 5694.375
$ws.Range("M77").Value = 118.8822500000006
$ws.Range("N77").Value = -14430.375
$ws.Range("H97").Value = 40311.58
$ws.Range("I97").Value = 53811.05
$ws.Range("J97").Value = 3670.1428
$ws.Range("K97").Value = 53811.05
$ws.Range("L97").Value = 3670.1428
$ws.Range("M97").Value = -53315.05
$ws.Range("N97").Value = -4662.1428
$ws.Range("H102").Value = 79348.766
$ws.Range("I102").Value = 127586.625
$ws.Range("J102").Value = 2168.2
$ws.Range("K102").Value = 127586.625
$ws.Range("L102").Value = 2168.2
$ws.Range("M102").Value = -125964.625
$ws.Range("N102").Value = -5412.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 712.0833
$ws.Range("I94").Value = 674.75
$ws.Range("J94").Value = 898.75
$ws.Range("K94").Value = 674.75
$ws.Range("L94").Value = 898.75
$ws.Range("M94").Value = -223.75
$ws.Range("N94").Value = -1800.75
$ws.Range("H99").Value = 2380
$ws.Range("I99").Value = 1726.6666
$ws.Range("J99").Value = 2520
$ws.Range("K99").Value = 1726.6666
$ws.Range("L99").Value = 2520
$ws.Range("M99").Value = -228.6666
$ws.Range("N99").Value = -5516
$ws.Range("H134").Value = 25996.738
$ws.Range("I134").Value = 35551.625
$ws.Range("K134").Value = 106654.875
$ws.Range("M134").Value = -104119.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 305
$ws.Range("I23").Value = 305
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 305
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -65
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 305
$ws.Range("I27").Value = 305
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 305
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -113
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 39406.5
$ws.Range("I31").Value = 642.3684
$ws.Range("J31").Value = 78170.63
$ws.Range("K31").Value = 642.3684
$ws.Range("L31").Value = 78170.63
$ws.Range("M31").Value = -347.3684
$ws.Range("N31").Value = -78760.63
$ws.Range("H34").Value = 39406.5
$ws.Range("I34").Value = 642.3684
$ws.Range("J34").Value = 78170.63
$ws.Range("K34").Value = 642.3684
$ws.Range("L34").Value = 78170.63
$ws.Range("M34").Value = -440.3684
$ws.Range("N34").Value = -78574.63
$ws.Range("H99").Value = 2864.6191
$ws.Range("I99").Value = 3123
$ws.Range("J99").Value = 2803.8235
$ws.Range("K99").Value = 3123
$ws.Range("L99").Value = 2803.8235
$ws.Range("M99").Value = -1625
$ws.Range("N99").Value = -5799.8235
$ws.Range("H126").Value = 2864.6191
$ws.Range("I126").Value = 3123
$ws.Range("J126").Value = 2803.8235
$ws.Range("K126").Value = 9369
$ws.Range("L126").Value = 8411.470499999999
$ws.Range("M126").Value = -6899
$ws.Range("N126").Value = -13351.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1320.1608
$ws.Range("J5").Value = 1333.2727
$ws.Range("L5").Value = 3999.8181
$ws.Range("N5").Value = -4223.8181
$ws.Range("H34").Value = 1058.84
$ws.Range("I34").Value = 518
$ws.Range("J34").Value = 1161.8572
$ws.Range("K34").Value = 1554
$ws.Range("L34").Value = 3485.5716
$ws.Range("M34").Value = -1470
$ws.Range("N34").Value = -3653.5716
$ws.Range("H58").Value = 1633.3334
$ws.Range("J58").Value = 1633.3334
$ws.Range("L58").Value = 4900.0002
$ws.Range("N58").Value = -5156.0002
$ws.Range("H135").Value = 1320.1608
$ws.Range("J135").Value = 1333.2727
$ws.Range("L135").Value = 11999.4543
$ws.Range("N135").Value = -17069.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2559.9644
$ws.Range("I102").Value = 1731.1052
$ws.Range("J102").Value = 4309.778
$ws.Range("K102").Value = 1731.1052
$ws.Range("L102").Value = 4309.778
$ws.Range("M102").Value = -109.1052
$ws.Range("N102").Value = -7553.778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 42504.68
$ws.Range("I40").Value = 68842.664
$ws.Range("J40").Value = 2997.7
$ws.Range("K40").Value = 68842.664
$ws.Range("L40").Value = 2997.7
$ws.Range("M40").Value = -68706.664
$ws.Range("N40").Value = -3269.7
$ws.Range("H61").Value = 2489.7273
$ws.Range("I61").Value = 2547.625
$ws.Range("J61").Value = 2335.3333
$ws.Range("K61").Value = 2547.625
$ws.Range("L61").Value = 2335.3333
$ws.Range("M61").Value = -2345.625
$ws.Range("N61").Value = -2739.3333
$ws.Range("H68").Value = 2040.9584
$ws.Range("I68").Value = 1823.5294
$ws.Range("J68").Value = 2569
$ws.Range("K68").Value = 1823.5294
$ws.Range("L68").Value = 2569
$ws.Range("M68").Value = -1074.5294
$ws.Range("N68").Value = -4067
$ws.Range("H71").Value = 2040.9584
$ws.Range("I71").Value = 1823.5294
$ws.Range("J71").Value = 2569
$ws.Range("K71").Value = 9117.646999999999
$ws.Range("L71").Value = 12845
$ws.Range("M71").Value = -5373.646999999999
$ws.Range("N71").Value = -20333
$ws.Range("H82").Value = 1994.8125
$ws.Range("I82").Value = 1641.8334
$ws.Range("J82").Value = 3053.75
$ws.Range("K82").Value = 1641.8334
$ws.Range("L82").Value = 3053.75
$ws.Range("M82").Value = -1280.8334
$ws.Range("N82").Value = -3775.75
$ws.Range("H85").Value = 1994.8125
$ws.Range("I85").Value = 1641.8334
$ws.Range("J85").Value = 3053.75
$ws.Range("K85").Value = 1641.8334
$ws.Range("L85").Value = 3053.75
$ws.Range("M85").Value = -393.8334
$ws.Range("N85").Value = -5549.75
$ws.Range("H93").Value = 2131.3076
$ws.Range("I93").Value = 2124.7778
$ws.Range("J93").Value = 2146
$ws.Range("K93").Value = 2124.7778
$ws.Range("L93").Value = 2146
$ws.Range("M93").Value = -876.7777999999998
$ws.Range("N93").Value = -4642
$ws.Range("H100").Value = 2169.1428
$ws.Range("I100").Value = 1901
$ws.Range("J100").Value = 2370.25
$ws.Range("K100").Value = 1901
$ws.Range("L100").Value = 2370.25
$ws.Range("M100").Value = -1360
$ws.Range("N100").Value = -3452.25
$ws.Range("H113").Value = 2489.7273
$ws.Range("I113").Value = 2547.625
$ws.Range("J113").Value = 2335.3333
$ws.Range("K113").Value = 2547.625
$ws.Range("L113").Value = 2335.3333
$ws.Range("M113").Value = -377.625
$ws.Range("N113").Value = -6675.3333
$ws.Range("H122").Value = 3159.682
$ws.Range("I122").Value = 2957.25
$ws.Range("J122").Value = 3699.5
$ws.Range("K122").Value = 8871.75
$ws.Range("L122").Value = 11098.5
$ws.Range("M122").Value = -6421.75
$ws.Range("N122").Value = -15998.5
$ws.Range("H132").Value = 3246
$ws.Range("I132").Value = 2986.2307
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 8958.6921
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -6428.6921
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 23182.5
$ws.Range("I56").Value = 3900
$ws.Range("K56").Value = 3900
$ws.Range("M56").Value = -3186
$ws.Range("H64").Value = 23000
$ws.Range("J64").Value = 23000
$ws.Range("L64").Value = 23000
$ws.Range("N64").Value = -23496
$ws.Range("H67").Value = 23000
$ws.Range("J67").Value = 23000
$ws.Range("L67").Value = 23000
$ws.Range("N67").Value = -24716
$ws.Range("H96").Value = 90910664
$ws.Range("I96").Value = 125001690
$ws.Range("J96").Value = 1268
$ws.Range("K96").Value = 125001690
$ws.Range("L96").Value = 1268
$ws.Range("M96").Value = -125000317
$ws.Range("N96").Value = -4014
$ws.Range("H141").Value = 60600
$ws.Range("J141").Value = 60600
$ws.Range("L141").Value = 60600
$ws.Range("N141").Value = -70960
